$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook used emoji markers in column A ("statut") to encode a status:
#   📘 -> ⚠️
#   📕 -> -3
#   📙 -> +3
#   📗 -> ✅
# Walk every used row in column A and replace the emoji value with its
# new text equivalent, forcing the "-3"/"+3" replacements to stay text
# (otherwise Excel would interpret them as negative/positive numbers).

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()

    if ($v -eq "📘") {
        $cell.Value = "⚠️"
    } elseif ($v -eq "📕") {
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
        $cell.Style = "Normal"
    } elseif ($v -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
        $cell.Style = "Normal"
    } elseif ($v -eq "📗") {
        $cell.Value = "✅"
    }
}
